# Update column F ("n" = response count) for specific rows in Sheet1.
# These rows correspond to the response category that picked up one
# additional respondent after excluding NA answers (report 2) / including
# one more response (report 3). Each listed cell increases by 1 relative
# to its original value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 25.0
$ws.Cells.Item(3, 6).Value = 25.0
$ws.Cells.Item(4, 6).Value = 25.0
$ws.Cells.Item(5, 6).Value = 25.0
$ws.Cells.Item(7, 6).Value = 17.0
$ws.Cells.Item(11, 6).Value = 17.0
$ws.Cells.Item(13, 6).Value = 13.0
$ws.Cells.Item(17, 6).Value = 19.0
$ws.Cells.Item(20, 6).Value = 11.0
$ws.Cells.Item(22, 6).Value = 14.0
$ws.Cells.Item(23, 6).Value = 23.0
$ws.Cells.Item(25, 6).Value = 25.0
$ws.Cells.Item(27, 6).Value = 5.0
$ws.Cells.Item(28, 6).Value = 26.0
$ws.Cells.Item(31, 6).Value = 8.0
$ws.Cells.Item(35, 6).Value = 26.0
$ws.Cells.Item(36, 6).Value = 25.0
$ws.Cells.Item(39, 6).Value = 8.0
$ws.Cells.Item(42, 6).Value = 16.0
$ws.Cells.Item(45, 6).Value = 7.0
$ws.Cells.Item(48, 6).Value = 7.0
$ws.Cells.Item(56, 6).Value = 2.0
$ws.Cells.Item(62, 6).Value = 7.0
$ws.Cells.Item(65, 6).Value = 6.0
$ws.Cells.Item(70, 6).Value = 6.0
$ws.Cells.Item(71, 6).Value = 10.0
$ws.Cells.Item(73, 6).Value = 26.0
$ws.Cells.Item(74, 6).Value = 12.0
$ws.Cells.Item(79, 6).Value = 10.0
$ws.Cells.Item(80, 6).Value = 15.0
$ws.Cells.Item(82, 6).Value = 9.0
$ws.Cells.Item(86, 6).Value = 5.0
$ws.Cells.Item(90, 6).Value = 4.0
$ws.Cells.Item(94, 6).Value = 12.0
$ws.Cells.Item(97, 6).Value = 26.0
$ws.Cells.Item(99, 6).Value = 21.0
$ws.Cells.Item(101, 6).Value = 22.0
$ws.Cells.Item(105, 6).Value = 6.0
$ws.Cells.Item(107, 6).Value = 26.0
$ws.Cells.Item(109, 6).Value = 12.0
$ws.Cells.Item(111, 6).Value = 11.0
$ws.Cells.Item(113, 6).Value = 12.0
$ws.Cells.Item(115, 6).Value = 10.0
$ws.Cells.Item(117, 6).Value = 21.0
$ws.Cells.Item(121, 6).Value = 21.0
$ws.Cells.Item(122, 6).Value = 24.0
$ws.Cells.Item(124, 6).Value = 26.0
$ws.Cells.Item(125, 6).Value = 26.0
$ws.Cells.Item(127, 6).Value = 16.0
$ws.Cells.Item(129, 6).Value = 9.0
$ws.Cells.Item(131, 6).Value = 11.0
$ws.Cells.Item(135, 6).Value = 15.0
$ws.Cells.Item(138, 6).Value = 18.0
$ws.Cells.Item(140, 6).Value = 26.0
$ws.Cells.Item(143, 6).Value = 18.0
$ws.Cells.Item(146, 6).Value = 6.0
$ws.Cells.Item(148, 6).Value = 10.0
$ws.Cells.Item(153, 6).Value = 4.0
$ws.Cells.Item(158, 6).Value = 2.0
$ws.Cells.Item(161, 6).Value = 11.0
$ws.Cells.Item(163, 6).Value = 13.0
$ws.Cells.Item(164, 6).Value = 17.0
$ws.Cells.Item(169, 6).Value = 11.0
$ws.Cells.Item(173, 6).Value = 7.0
$ws.Cells.Item(178, 6).Value = 4.0
$ws.Cells.Item(179, 6).Value = 25.0
$ws.Cells.Item(180, 6).Value = 22.0
$ws.Cells.Item(184, 6).Value = 5.0
$ws.Cells.Item(189, 6).Value = 4.0
$ws.Cells.Item(193, 6).Value = 16.0
$ws.Cells.Item(197, 6).Value = 8.0
$ws.Cells.Item(201, 6).Value = 8.0
$ws.Cells.Item(203, 6).Value = 18.0
$ws.Cells.Item(205, 6).Value = 14.0
$ws.Cells.Item(209, 6).Value = 8.0
$ws.Cells.Item(212, 6).Value = 13.0
$ws.Cells.Item(217, 6).Value = 4.0
$ws.Cells.Item(221, 6).Value = 11.0
$ws.Cells.Item(224, 6).Value = 6.0
$ws.Cells.Item(226, 6).Value = 25.0
$ws.Cells.Item(227, 6).Value = 25.0
$ws.Cells.Item(229, 6).Value = 12.0
$ws.Cells.Item(231, 6).Value = 7.0
$ws.Cells.Item(232, 6).Value = 26.0
$ws.Cells.Item(233, 6).Value = 26.0
$ws.Cells.Item(235, 6).Value = 15.0
$ws.Cells.Item(237, 6).Value = 12.0
$ws.Cells.Item(241, 6).Value = 3.0
$ws.Cells.Item(243, 6).Value = 25.0
$ws.Cells.Item(245, 6).Value = 11.0
$ws.Cells.Item(246, 6).Value = 25.0
$ws.Cells.Item(248, 6).Value = 17.0
$ws.Cells.Item(250, 6).Value = 15.0
$ws.Cells.Item(254, 6).Value = 11.0
